# The deck's single applied design ("Integral") is switched to the
# built-in "Office Theme" colour palette -- i.e. the presentation's
# theme colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) is
# replaced with the standard Office theme colours.
#
# PowerPoint's ThemeColorScheme exposes exactly those twelve slots, in
# this order:
#   1 dk1  2 lt1  3 dk2  4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#   11 hlink  12 folHlink

$p = $ppt.ActivePresentation

# Target palette = the stock "Office Theme" colour scheme.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # COM RGB() packing: R + G*256 + B*65536
    $themeColors.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
